$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# Simple "want to go" count (F column) updates
$ws1.Range("F3").Value  = 26981
$ws1.Range("F5").Value  = 259
$ws1.Range("F6").Value  = 628
$ws1.Range("F10").Value = 370
$ws1.Range("F11").Value = 462
$ws1.Range("F12").Value = 194
$ws1.Range("F14").Value = 309
$ws1.Range("F15").Value = 92
$ws1.Range("F16").Value = 457
$ws1.Range("F18").Value = 1587
$ws1.Range("F19").Value = 233
$ws1.Range("F20").Value = 119
$ws1.Range("F21").Value = 449

# Insert a brand-new event row at row 22 ("第十届萌物语动漫嘉年华"), pushing the
# existing rows 22-23 down to 23-24.
$ws1.Rows.Item(22).Insert()

# Copy formatting (border/bold/alignment) of the index column from the row
# directly above so the new row matches the existing look of column A.
$ws1.Range("A21").Copy()
$ws1.Range("A22").PasteSpecial(-4122)

$ws1.Range("A22").Value = 21
$ws1.Range("B22").NumberFormat = "@"
$ws1.Range("B22").Value = "2024-10-07"
$ws1.Range("B22").NumberFormat = "General"
$ws1.Range("C22").Value = "广州·第十届萌物语动漫嘉年华"
$ws1.Range("D22").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws1.Range("E22").Value = "2024.10.07 10:00-10.07 17:00"
$ws1.Range("F22").Value = 2
$ws1.Range("G22").Value = 9.9
$ws1.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=91162"
$ws1.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202408/9m6CSSzf1723739254235.jpeg"

# Fix up the index numbers + "want to go" count on the rows that shifted down
$ws1.Range("A23").Value = 22
$ws1.Range("F23").Value = 108
$ws1.Range("A24").Value = 23

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4517

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5138
$ws3.Range("F3").Value = 259

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) -- same underlying data, duplicated layout
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 5138
$ws4.Range("F4").Value  = 259
$ws4.Range("F5").Value  = 26981
$ws4.Range("F7").Value  = 4517
$ws4.Range("F8").Value  = 259
$ws4.Range("F10").Value = 628
$ws4.Range("F23").Value = 370
$ws4.Range("F24").Value = 462
$ws4.Range("F25").Value = 194
$ws4.Range("F28").Value = 309
$ws4.Range("F29").Value = 92
$ws4.Range("F32").Value = 457
$ws4.Range("F35").Value = 1587
$ws4.Range("F36").Value = 233
$ws4.Range("F38").Value = 119
$ws4.Range("F39").Value = 449

# Insert the same new event row at row 40, pushing rows 40-47 down to 41-48.
$ws4.Rows.Item(40).Insert()

$ws4.Range("A39").Copy()
$ws4.Range("A40").PasteSpecial(-4122)

$ws4.Range("A40").Value = 39
$ws4.Range("B40").NumberFormat = "@"
$ws4.Range("B40").Value = "2024-10-07"
$ws4.Range("B40").NumberFormat = "General"
$ws4.Range("C40").Value = "广州·第十届萌物语动漫嘉年华"
$ws4.Range("D40").Value = "奥体南路12号 优托邦(奥体旗舰店)"
$ws4.Range("E40").Value = "2024.10.07 10:00-10.07 17:00"
$ws4.Range("F40").Value = 2
$ws4.Range("G40").Value = 9.9
$ws4.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=91162"
$ws4.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202408/9m6CSSzf1723739254235.jpeg"

$ws4.Range("A41").Value = 40
$ws4.Range("F41").Value = 108
$ws4.Range("A42").Value = 41
$ws4.Range("A43").Value = 42
$ws4.Range("A44").Value = 43
$ws4.Range("A45").Value = 44
$ws4.Range("A46").Value = 45
$ws4.Range("A47").Value = 46
$ws4.Range("A48").Value = 47
